# update yellow keys when player change
$wb = $excel.ActiveWorkbook

# --- VMIX sheet: row 2 (final section summary) ---
$ws = $wb.Worksheets.Item("VMIX")

# Flags for finalists 1 and 2 swapped
$ws.Range("CZ2").Value = "C:\TRIAL_2021\VMIX\MATERIAL\BANDERES\esp.png"
$ws.Range("DA2").Value = "C:\TRIAL_2021\VMIX\MATERIAL\BANDERES\fra.png"

# Countries for finalists 1 and 2 swapped
$ws.Range("DF2").Value = "ESP"
$ws.Range("DG2").Value = "FRA"

# Player names for finalists 1 and 2 swapped
$ws.Range("DL2").Value = "ALEJANDRO MO"
$ws.Range("DM2").Value = "VINCENT H"

# Points for finalists 1 and 2
$ws.Range("DR2").Value = 60
$ws.Range("DS2").Value = 20

# Section 1 scores
$ws.Range("DY2").Value = 0
$ws.Range("EP2").Value = 30
$ws.Range("EQ2").Value = 20

# Abbreviations for finalists 1 and 2 swapped
$ws.Range("FB2").Value = "MON"
$ws.Range("FC2").Value = "HER"

# Current section info
$ws.Range("FH2").Value = "SECTION 1"
$ws.Range("FI2").Value = 30
$ws.Range("FJ2").Value = 10
$ws.Range("FK2").Value = 10
$ws.Range("FL2").Value = 10
$ws.Range("FQ2").Value = 1
$ws.Range("FR2").Value = 7
$ws.Range("FT2").Value = "ALEJANDRO MO"
$ws.Range("FV2").Value = "SECTION 4"
$ws.Range("FX2").Value = 20

# --- TRIAL sheet: rows 18-19 (players section scores) ---
$ws2 = $wb.Worksheets.Item("TRIAL")

$ws2.Range("H18").Value = 0
$ws2.Range("K18").Value = 20
$ws2.Range("M18").Value = 20
$ws2.Range("Q18").Value = 0
$ws2.Range("R18").Value = 1
$ws2.Range("U18").Value = 1

$ws2.Range("K19").Value = 30
$ws2.Range("M19").Value = 60
$ws2.Range("Q19").Value = 2
$ws2.Range("T19").Value = 3
$ws2.Range("U19").Value = 0
